# Aggiornamento dati fino al 13/05 (aggiunta righe 252-255)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

$lastRow = 251

$data = @(
    @(44326, 5, 21, 131.4965560425798),
    @(44327, 4, 23, 144.0200375704446),
    @(44328, 0, 23, 144.0200375704446),
    @(44329, 3, 20, 125.2348152786475)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $lastRow + 1 + $i
    $values = $data[$i]

    # Copy the formatting (style) used by column A of the previous row
    # so the new date cell keeps the same style (s="2": centered, bold,
    # bordered, custom date number format) as every other row above.
    $ws.Range("A" + ($row - 1)).Copy()
    $ws.Range("A" + $row).PasteSpecial($xlPasteFormats)
    $excel.CutCopyMode = 0

    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
